{"js": "// Edit: split the paragraph that reads \"ghjhh\" into three runs:\n//   \"G\"  +  \"hjhh\"  +  \"n   hgjhkjh\"\n// (net visible text becomes \"Ghjhhn   hgjhkjh\"), matching the target\n// OOXML diff which turns a single <w:r> into three sibling <w:r> elements.\n//\n// A plain `range.insertText(..., \"Replace\")` (or setting `.text`) would\n// collapse back down to a single run because the replacement text shares\n// the same run formatting as its neighbours, so the run split would be\n// lost. To faithfully reproduce the exact run boundaries from the diff we\n// splice in literal OOXML for the paragraph's run content via\n// `Range.insertOoxml`, which Word does not auto-merge.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph whose full text is exactly \"ghjhh\" (the pre-edit\n// state from the diff) so the script is resilient to the paragraph's\n// position in the document.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"ghjhh\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the target paragraph containing 'ghjhh'.\");\n}\n\n// `Range.insertOoxml` requires a flat-OPC `<pkg:package>` wrapper around\n// the part's XML (Office.js validates this before sending it through).\n// Only the run-level content of the paragraph is replaced; the paragraph\n// itself (and its mark/properties) is kept because `getRange()` on a\n// paragraph addresses its content, not the surrounding structure.\nconst flatOpcXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t>G</w:t></w:r>' +\n  '<w:r><w:t>hjhh</w:t></w:r>' +\n  '<w:r><w:t>n   hgjhkjh</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ntarget.getRange().insertOoxml(flatOpcXml, \"Replace\");\nawait context.sync();\n", "ps1": "# Edit: split the paragraph that reads \"ghjhh\" into three runs:\n#   \"G\"  +  \"hjhh\"  +  \"n   hgjhkjh\"\n# (net visible text becomes \"Ghjhhn   hgjhkjh\"), matching the target\n# OOXML diff which turns a single <w:r> into three sibling <w:r> elements.\n#\n# Setting $range.Text (or a Find/Replace) would just rewrite the single\n# run's text and would not reproduce the three separate <w:r> elements the\n# diff calls for, because Word merges same-formatted adjacent text back\n# into one run. Splicing literal WordprocessingML via Range.InsertXML\n# preserves the exact run boundaries instead.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    # Paragraph.Range.Text includes the trailing paragraph-mark character,\n    # so trim it before comparing against the plain text from the diff.\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t -eq \"ghjhh\") {\n        $target = $p\n        break\n    }\n}\nif ($null -eq $target) {\n    throw \"Could not find the target paragraph containing 'ghjhh'.\"\n}\n\n$runXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:r><w:t>G</w:t></w:r>' +\n          '<w:r><w:t>hjhh</w:t></w:r>' +\n          '<w:r><w:t>n   hgjhkjh</w:t></w:r>' +\n          '</w:p>'\n\n$target.Range.InsertXML($runXml)\n"}
